$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracker")
$lo = $ws.ListObjects.Item(1)

# --- Row 4: naive XGB test run (standardized scaling only) ---
$ws.Range("A4").Value = Get-Date -Year 2023 -Month 2 -Day 14 -Hour 21 -Minute 51 -Second 0
$ws.Range("B4").Value = "Test"
$ws.Range("D4").Value = "XGB"
$ws.Range("E4").Value = "scaling"
$ws.Range("H4").Value = 0.989
$ws.Range("I4").Value = 0.857
$ws.Range("K4").Value = "Maria"

# --- Row 5: naive XGB test run with first hyperparameter set ---
$ws.Range("A5").Value = Get-Date -Year 2023 -Month 2 -Day 14 -Hour 22 -Minute 4 -Second 0
$ws.Range("B5").Value = "Test"
$ws.Range("D5").Value = "XGB"
$ws.Range("E5").Value = "scaling"
$ws.Range("G5").Value = '{"model__colsample_bylevel": 0.3, "model__max_depth": 7}'
$ws.Range("H5").Value = 0.996
$ws.Range("I5").Value = 0.856
$ws.Range("K5").Value = "Maria"

# --- Row 6: submission of that naive XGB model ---
$ws.Range("A6").Value = Get-Date -Year 2023 -Month 2 -Day 14 -Hour 22 -Minute 5 -Second 0
$ws.Range("B6").Value = "Submission"
$ws.Range("C6").Value = "20230214_2204_naive_xgb.csv"
$ws.Range("D6").Value = "XGB"
$ws.Range("E6").Value = "scaling"
$ws.Range("G6").Value = '{"model__colsample_bylevel": 0.3, "model__max_depth": 7}'
$ws.Range("H6").Value = 0.991
$ws.Range("J6").Value = 0.556
$ws.Range("K6").Value = "Maria"

# --- Row 7: naive XGB test run with second hyperparameter set ---
$ws.Range("A7").Value = Get-Date -Year 2023 -Month 2 -Day 14 -Hour 22 -Minute 20 -Second 0
$ws.Range("B7").Value = "Test"
$ws.Range("D7").Value = "XGB"
$ws.Range("E7").Value = "scaling"
$ws.Range("G7").Value = '{"model__colsample_bylevel": 0.2, "model__min_child_weight": 50, "model__max_depth": 5}'
$ws.Range("H7").Value = 0.865
$ws.Range("I7").Value = 0.801
$ws.Range("K7").Value = "Maria"

# --- Row 8: submission of that second naive XGB model ---
$ws.Range("A8").Value = Get-Date -Year 2023 -Month 2 -Day 14 -Hour 22 -Minute 21 -Second 0
$ws.Range("B8").Value = "Submission"
$ws.Range("D8").Value = "XGB"
$ws.Range("E8").Value = "scaling"
$ws.Range("G8").Value = '{"model__colsample_bylevel": 0.2, "model__min_child_weight": 50, "model__max_depth": 5}'
$ws.Range("H8").Value = 0.863
$ws.Range("J8").Value = 0.50765
$ws.Range("J8").NumberFormat = "0.000"
$ws.Range("K8").Value = "Maria"

# --- Grow the tracker table by one blank row (keeps formatting consistent) ---
$newRow = $lo.ListRows.Add()
$ws.Range("A3:C3").Copy()
$ws.Range("A31:C31").PasteSpecial(-4122)
$ws.Range("H3:I3").Copy()
$ws.Range("H31:I31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Restore view state (zoom + selection) ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 108
$ws.Range("G5").Select()
